$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = -1707
$ws.Range("H33").Value = 613.8570999999999
$ws.Range("I33").Value = 135.53847
$ws.Range("J33").Value = 1391.125
$ws.Range("K33").Value = 135.53847
$ws.Range("L33").Value = 1391.125
$ws.Range("M33").Value = 93.46153000000001
$ws.Range("N33").Value = -1849.125
$ws.Range("H62").Value = 1908.2
$ws.Range("I62").Value = 1885.375
$ws.Range("J62").Value = 1999.5
$ws.Range("K62").Value = 1885.375
$ws.Range("L62").Value = 1999.5
$ws.Range("M62").Value = -1261.375
$ws.Range("N62").Value = -3247.5
$ws.Range("H65").Value = 1908.2
$ws.Range("I65").Value = 1885.375
$ws.Range("J65").Value = 1999.5
$ws.Range("K65").Value = 9426.875
$ws.Range("L65").Value = 9997.5
$ws.Range("M65").Value = -6306.875
$ws.Range("N65").Value = -16237.5
$ws.Range("H76").Value = 3323.9443
$ws.Range("I76").Value = 2987.9285
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 2987.9285
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -2672.9285
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 3323.9443
$ws.Range("I79").Value = 2987.9285
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 2987.9285
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -1895.9285
$ws.Range("N79").Value = -6684
$ws.Range("H103").Value = 878.55554
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 878.55554
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = $null
$ws.Range("M103").Value = 2635.66662
$ws.Range("N103").Value = -3807.66662
$ws.Range("H132").Value = 5251.0347
$ws.Range("I132").Value = 7073.8
$ws.Range("J132").Value = 1200.4445
$ws.Range("K132").Value = 21221.4
$ws.Range("L132").Value = 3601.3335
$ws.Range("M132").Value = -18691.4
$ws.Range("N132").Value = -8661.333500000001
$ws.Range("H137").Value = 1371.1698
$ws.Range("I137").Value = 954.1429000000001
$ws.Range("K137").Value = 2862.4287
$ws.Range("M137").Value = -312.4287000000004
$ws.Range("H138").Value = 3333.3262
$ws.Range("I138").Value = 4051.7144
$ws.Range("J138").Value = 3204.3845
$ws.Range("K138").Value = 12155.1432
$ws.Range("L138").Value = 9613.1535
$ws.Range("M138").Value = -7015.143199999999
$ws.Range("N138").Value = -19893.1535

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25778.578
$ws.Range("I32").Value = 4807.574
$ws.Range("J32").Value = 92392.35000000001
$ws.Range("K32").Value = 4807.574
$ws.Range("L32").Value = 92392.35000000001
$ws.Range("M32").Value = -4520.574
$ws.Range("N32").Value = -92966.35000000001
$ws.Range("H44").Value = 12791.25
$ws.Range("J44").Value = 12791.25
$ws.Range("L44").Value = 12791.25
$ws.Range("N44").Value = -13767.25
$ws.Range("H45").Value = 1520.125
$ws.Range("I45").Value = 970.5
$ws.Range("J45").Value = 2069.75
$ws.Range("K45").Value = 970.5
$ws.Range("L45").Value = 2069.75
$ws.Range("M45").Value = -593.5
$ws.Range("N45").Value = -2823.75
$ws.Range("H102").Value = 33130.906
$ws.Range("I102").Value = 47200.547
$ws.Range("K102").Value = 47200.547
$ws.Range("M102").Value = -45578.547
$ws.Range("H132").Value = 4875.9585
$ws.Range("I132").Value = 5379.579
$ws.Range("K132").Value = 16138.737
$ws.Range("M132").Value = -13608.737

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 22251330
$ws.Range("I107").Value = 30340716
$ws.Range("K107").Value = 30340716
$ws.Range("M107").Value = -30338796
$ws.Range("H134").Value = 1856.871
$ws.Range("I134").Value = 1728.7037
$ws.Range("J134").Value = 2722
$ws.Range("K134").Value = 5186.1111
$ws.Range("L134").Value = 8166
$ws.Range("M134").Value = -2651.1111
$ws.Range("N134").Value = -13236

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = $null
$ws.Range("N54").Value = 0
$ws.Range("H107").Value = 822.2941
$ws.Range("I107").Value = 948.1111
$ws.Range("J107").Value = 680.75
$ws.Range("K107").Value = 948.1111
$ws.Range("L107").Value = 680.75
$ws.Range("M107").Value = 971.8889
$ws.Range("N107").Value = -4520.75
$ws.Range("H122").Value = 2685.5334
$ws.Range("I122").Value = 3005.8333
$ws.Range("K122").Value = 9017.499899999999
$ws.Range("M122").Value = -6567.499899999999
$ws.Range("H132").Value = 1743.7576
$ws.Range("I132").Value = 907.64
$ws.Range("K132").Value = 2722.92
$ws.Range("M132").Value = -192.9200000000001
$ws.Range("H134").Value = 1123.9487
$ws.Range("I134").Value = 862.3823
$ws.Range("J134").Value = 2902.6
$ws.Range("K134").Value = 2587.1469
$ws.Range("L134").Value = 8707.799999999999
$ws.Range("M134").Value = -52.14689999999973
$ws.Range("N134").Value = -13777.8

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8806.360000000001
$ws.Range("I5").Value = 620.875
$ws.Range("J5").Value = 23358.334
$ws.Range("K5").Value = 1862.625
$ws.Range("L5").Value = 70075.00199999999
$ws.Range("M5").Value = -1750.625
$ws.Range("N5").Value = -70299.00199999999
$ws.Range("H68").Value = 17440.92
$ws.Range("I68").Value = 953.95
$ws.Range("J68").Value = 25291.857
$ws.Range("K68").Value = 2861.85
$ws.Range("L68").Value = 75875.571
$ws.Range("M68").Value = -2050.85
$ws.Range("N68").Value = -77497.571
$ws.Range("H71").Value = 17440.92
$ws.Range("I71").Value = 953.95
$ws.Range("J71").Value = 25291.857
$ws.Range("K71").Value = 8585.550000000001
$ws.Range("L71").Value = 227626.713
$ws.Range("M71").Value = -4529.550000000001
$ws.Range("N71").Value = -235738.713
$ws.Range("H109").Value = 3837.9473
$ws.Range("I109").Value = 1393.4
$ws.Range("J109").Value = 4711
$ws.Range("K109").Value = 4180.200000000001
$ws.Range("L109").Value = 14133
$ws.Range("M109").Value = -3140.200000000001
$ws.Range("N109").Value = -16213
$ws.Range("H122").Value = 764.5
$ws.Range("J122").Value = 1249.5
$ws.Range("L122").Value = 11245.5
$ws.Range("N122").Value = -16145.5
$ws.Range("H135").Value = 8806.360000000001
$ws.Range("I135").Value = 620.875
$ws.Range("J135").Value = 23358.334
$ws.Range("K135").Value = 5587.875
$ws.Range("L135").Value = 210225.006
$ws.Range("M135").Value = -3052.875
$ws.Range("N135").Value = -215295.006

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2479.28
$ws.Range("I132").Value = 2363.1428
$ws.Range("J132").Value = 2627.0908
$ws.Range("K132").Value = 7089.428400000001
$ws.Range("L132").Value = 7881.2724
$ws.Range("M132").Value = -4559.428400000001
$ws.Range("N132").Value = -12941.2724

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = $null
$ws.Range("N34").Value = 0

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 19076
$ws.Range("J25").Value = 19076
$ws.Range("L25").Value = 19076
$ws.Range("N25").Value = -19662
$ws.Range("H32").Value = 6800
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 16400
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 16400
$ws.Range("M32").Value = -1683
$ws.Range("N32").Value = -17034
$ws.Range("H34").Value = 16800
$ws.Range("J34").Value = 16800
$ws.Range("L34").Value = 16800
$ws.Range("N34").Value = -17206
$ws.Range("H119").Value = 41992.5
$ws.Range("J119").Value = 41992.5
$ws.Range("L119").Value = 41992.5
$ws.Range("N119").Value = -51668.5
$ws.Range("H136").Value = 1233
$ws.Range("I136").Value = 627.38464
$ws.Range("J136").Value = 2357.7144
$ws.Range("K136").Value = 1882.15392
$ws.Range("L136").Value = 7073.1432
$ws.Range("M136").Value = 667.84608
$ws.Range("N136").Value = -12173.1432
